$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 9 (2oKzt7dm / COLOMBIA - PRIMERA A / Pereira vs Deportes Tolima)
# This shifts all subsequent rows up by one, matching the diff which removes
# the last row (24) and shrinks the used range from A1:AS24 to A1:AS23.
$ws.Rows("9:9").Delete()
